# Xây dựng hệ thống
#
# After the "Giai đoạn 4:..." paragraph (the last body paragraph, which also
# carries the stray "_GoBack" bookmark left over from the last edit position),
# insert four new paragraphs that begin section 7 ("Phân tích thiết kê") of
# the document, and move the "_GoBack" bookmark onto the final one of those
# new paragraphs ("7.2.Giao diện"), matching where Word would have left the
# caret after typing the new content.

$d = $word.ActiveDocument

# 1) The "_GoBack" bookmark currently sits at the end of the "Giai đoạn 4"
#    paragraph. Remove it there -- it will be re-created on the last new
#    paragraph below.
try {
    $oldMark = $d.Bookmarks.Item("_GoBack")
    $oldMark.Delete()
} catch {
    # No pre-existing "_GoBack" bookmark -- nothing to remove.
}

# 2) Locate the end of the "Giai đoạn 4" text -- this is where the new
#    paragraphs get inserted.
$find = $d.Content
$found = $find.Find.Execute(
    "Giai đoạn 4:Testing và chỉnh sửa .Và đưa vào vận hành",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the 'Giai đoạn 4' paragraph to anchor the new content."
}

# Re-seat the insertion point into a fresh Range object (positioned right
# after the matched text) so it isn't still tied to the Find match span.
$insertAt = $find.End
$target = $d.Range($insertAt, $insertAt)

$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$heading = '<w:p ' + $w + '><w:pPr><w:rPr><w:b/><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr><w:t>7</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr><w:t>&gt;Phân tích thiết kê</w:t></w:r></w:p>'

$p71 = '<w:p ' + $w + '><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>7.1.Mô hình tích hợp phần cứng phần mềm</w:t></w:r></w:p>'

$pHw = '<w:p ' + $w + '><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>-Phần cứng rất quan trọng.Nó quyết định đến tốc độ, hiệu xuất, tính tiện lợi và an toàn của hệ thống.</w:t></w:r></w:p>'

$p72 = '<w:p ' + $w + '><w:pPr><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:lastRenderedPageBreak/><w:t>7.2.Giao diện</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

[void]$target.InsertXml($heading + $p71 + $pHw + $p72)
